$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data row (row 2) with new values
$ws.Range("A2").Value = "ebenito"
$ws.Range("C2").Value = "'AA22025GL1MH"
$ws.Range("E2").Value = "PASSED"
$ws.Range("F2").Value = "TT232005HPMK 11:2"
$ws.Range("G2").Value = "19 jul. 2023, 11:27:20"

# Autofit column E since it now has content that affects its width
$ws.Columns.Item(5).EntireColumn.AutoFit()

# Update the active selection to B2
$ws.Range("B2").Select()
